# PP2 - Create standard workflows
# Applies changes to the "openml_100" worksheet:
#  - Set D column to "Run" for several rows (new task runs added)
#  - Change some "ensemble_snapshot" Configuration values to "standard"
#  - Change row 33 from "progressing" to "complete" and clear its Run flag
#  - Update the active cell selection on the sheet

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openml_100")

# Add "Run" marker (column D) to rows that are newly scheduled
$runRows = @(4, 9, 20, 53, 71, 76, 79)
foreach ($r in $runRows) {
    $ws.Cells.Item($r, 4).Value = "Run"
}

# Switch Configuration (column B) from "ensemble_snapshot" to "standard"
$standardRows = @(7, 16, 61, 62)
foreach ($r in $standardRows) {
    $ws.Cells.Item($r, 2).Value = "standard"
}

# Row 33 finished running: mark Status as "complete" and clear the Run flag
$ws.Cells.Item(33, 3).Value = "complete"
$ws.Cells.Item(33, 4).ClearContents()

# Update the active selection to reflect where the user left off
$ws.Range("C80").Select()
